$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 702
$ws.Range("I2").Value = 1776
$ws.Range("J2").Value = 7534
$ws.Range("K2").Value = 34
$ws.Range("L2").Value = 2125
$ws.Range("M2").Value = 112
$ws.Range("N2").Value = 1300
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 107
$ws.Range("S2").Value = 815
$ws.Range("T2").Value = 1371
$ws.Range("U2").Value = 110
$ws.Range("V2").Value = 11633
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 11635
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 191
$ws.Range("AA2").Value = 77
